# Remove the blank "Sheet1" tab, and set the Miscellaneous sector
# correction row to explicit zeros (instead of blank) on CO, NOX, SO2, VOC.
# Also touch up the shared-formula break on NH3-Org_and_Adj!G42.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the empty "Sheet1" tab -------------------------------------
$excel.DisplayAlerts = $false
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()
$excel.DisplayAlerts = $true

# --- 2. Miscellaneous sector correction row: blank -> 0 -------------------
# CO, NOX, SO2, VOC all have a "Miscellaneous" correction row that was left
# blank; fill it in with explicit zeros across the historical-year columns.

$targets = @(
    @{ Sheet = "CO";  Row = 18; LastCol = "AD" },
    @{ Sheet = "NOX"; Row = 18; LastCol = "AD" },
    @{ Sheet = "SO2"; Row = 18; LastCol = "AD" },
    @{ Sheet = "VOC"; Row = 17; LastCol = "AD" }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    $rng = $ws.Range("B$($t.Row):$($t.LastCol)$($t.Row)")
    $rng.Value = 0
}

# --- 3. NH3-Org_and_Adj!G42: re-enter the formula explicitly ---------------
$nh3org = $wb.Worksheets.Item("NH3-Org_and_Adj")
$nh3org.Range("G42").Formula = "=G20*(`$L`$37/`$L`$36)"

# --- 4. Make NH3-Org_and_Adj the active sheet/tab --------------------------
$nh3org.Select()
$nh3org.Range("B18:B21").Select()
$nh3org.Application.ActiveWindow.ScrollRow = 7
